# NVE data download and plotting
# Splits the old "NVE_datadownload_01-2022" sheet into three purpose-specific
# sheets (discharge, stage-only, suspended-sediment) and tweaks a couple of
# unrelated view/format bits that were touched in the same session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Hydrol sheet: widen column F a bit and leave the cursor on the last
#    (whole) data row.
# ---------------------------------------------------------------------------
$wsHydrol = $wb.Worksheets.Item("Hydrol")
$wsHydrol.Columns.Item(6).ColumnWidth = 15
$wsHydrol.Rows.Item(14).Select()

# ---------------------------------------------------------------------------
# 2. Rework the NVE sheet: rename to "nve_discharge", retype its content
#    (adding the Vækkava station and moving Julelv to the bottom), and add
#    two sibling sheets for the other data types.
# ---------------------------------------------------------------------------
$wsDischarge = $wb.Worksheets.Item("NVE_datadownload_01-2022")
$wsDischarge.Name = "nve_discharge"

$wsDischarge.Cells.Clear()
$dischargeData = @(
    @("stasjonID", "stasjonNavn"),
    @("234.14.0", "Cærrogæsjokka"),
    @("234.2.0", "Jiesjokka"),
    @("234.5.0", "Njalmigoaika"),
    @("234.1.0", "Polmak"),
    @("234.18.0", "Polmak nye"),
    @("234.4.0", "Smalfjord"),
    @("234.13.0", "Veahkkava, Iesjokka"),
    @("234.6.0", "Vækkava"),
    @("234.16.0", "Julelv")
)
for ($i = 0; $i -lt $dischargeData.Length; $i++) {
    $row = $i + 1
    $wsDischarge.Cells.Item($row, 1).Value = $dischargeData[$i][0]
    $wsDischarge.Cells.Item($row, 2).Value = $dischargeData[$i][1]
}
$wsDischarge.Range("A10:B10").Style = "Normal 3"

# New sheet: stage-only stations
$wsStage = $wb.Worksheets.Add($null, $wsDischarge)
$wsStage.Name = "nve_stage_only"
$stageData = @(
    @("stasjonID", "stasjonNavn"),
    @("234.10.0", "Karasjok"),
    @("234.12.0", "Nedre Levajok"),
    @("234.3.0", "Jiesjavrre")
)
for ($i = 0; $i -lt $stageData.Length; $i++) {
    $row = $i + 1
    $wsStage.Cells.Item($row, 1).Value = $stageData[$i][0]
    $wsStage.Cells.Item($row, 2).Value = $stageData[$i][1]
}
$wsStage.Range("L41").Select()

# New sheet: suspended sediment stations
$wsSS = $wb.Worksheets.Add($null, $wsStage)
$wsSS.Name = "nve_SS_1208-1200"
$ssData = @(
    @("stasjonID", "stasjonNavn"),
    @("234.1.0", "Polmak"),
    @("234.18.0", "Polmak nye"),
    @("234.32.0", "Tana v/Storfossen"),
    @("234.16.0", "Julelv")
)
for ($i = 0; $i -lt $ssData.Length; $i++) {
    $row = $i + 1
    $wsSS.Cells.Item($row, 1).Value = $ssData[$i][0]
    $wsSS.Cells.Item($row, 2).Value = $ssData[$i][1]
}
$wsSS.Range("A5:B5").Style = "Normal 3"
$wsSS.Range("R36").Select()

# ---------------------------------------------------------------------------
# 3. Make nve_discharge the active tab (it was previously Summary_forJLG_01-2022)
# ---------------------------------------------------------------------------
$wsDischarge.Activate()
$wsDischarge.Range("N28").Select()

Write-Output "NVE data download and plotting edit applied"
